# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-14 09:24:04
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists
# the users who touched each session record. For a handful of rows the
# synced/upstream data now lists the non-"System" contributor(s) first and
# moves "System" (and any trailing duplicate "system" entry) after it, i.e.
# the last comma-separated entry is rotated to the front of the list.
#
# Only the specific rows below changed between the two syncs; every other
# "System, ..." row on the sheet is untouched by this sync, so we update
# cell-by-cell rather than rewriting the whole column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @(
    @{ Row = 2;   Value = "backup@backdoor.com, System, system" }
    @{ Row = 3;   Value = "dnasr281@gmail.com, System" }
    @{ Row = 4;   Value = "backup@backdoor.com, System" }
    @{ Row = 5;   Value = "backup@backdoor.com, System" }
    @{ Row = 6;   Value = "dnasr281@gmail.com, System" }
    @{ Row = 7;   Value = "admin@admin.com, System" }
    @{ Row = 8;   Value = "backup@backdoor.com, System" }
    @{ Row = 28;  Value = "backup@backdoor.com, System, system" }
    @{ Row = 29;  Value = "dnasr281@gmail.com, System" }
    @{ Row = 30;  Value = "backup@backdoor.com, System" }
    @{ Row = 31;  Value = "backup@backdoor.com, System" }
    @{ Row = 32;  Value = "dnasr281@gmail.com, System" }
    @{ Row = 33;  Value = "admin@admin.com, System" }
    @{ Row = 34;  Value = "backup@backdoor.com, System" }
    @{ Row = 54;  Value = "backup@backdoor.com, System, system" }
    @{ Row = 55;  Value = "dnasr281@gmail.com, System" }
    @{ Row = 56;  Value = "backup@backdoor.com, System" }
    @{ Row = 57;  Value = "backup@backdoor.com, System" }
    @{ Row = 58;  Value = "dnasr281@gmail.com, System" }
    @{ Row = 59;  Value = "admin@admin.com, System" }
    @{ Row = 60;  Value = "backup@backdoor.com, System" }
    @{ Row = 80;  Value = "backup@backdoor.com, System" }
    @{ Row = 81;  Value = "backup@backdoor.com, System" }
    @{ Row = 82;  Value = "backup@backdoor.com, System" }
    @{ Row = 106; Value = "backup@backdoor.com, System" }
    @{ Row = 107; Value = "backup@backdoor.com, System" }
    @{ Row = 108; Value = "backup@backdoor.com, System" }
    @{ Row = 132; Value = "backup@backdoor.com, System" }
    @{ Row = 133; Value = "backup@backdoor.com, System" }
    @{ Row = 134; Value = "backup@backdoor.com, System" }
)

foreach ($update in $updates) {
    $cell = $ws.Cells.Item($update.Row, 7)  # Column G = "Recorded By"
    $cell.Value = $update.Value
}
